$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 245, pushing existing rows 245-260 down to 246-261.
$ws.Rows.Item(245).Insert()

# Populate the newly inserted row 245 with the new weekly record.
$ws.Range("A245").Value = 10
$ws.Range("B245").Value = "Vega Modelo de Temuco"
$ws.Range("C245").Value = "La Araucanía"
$ws.Range("D245").Value = 44585
$ws.Range("E245").Value = 9
$ws.Range("F245").Value = 100112009
$ws.Range("G245").Value = "Acelga"
$ws.Range("H245").Value = "Sin especificar"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 65
$ws.Range("K245").Value = 7000
$ws.Range("L245").Value = 7000
$ws.Range("M245").Value = 7000
$ws.Range("N245").Value = "$/docena de atados (12 kilos)"
$ws.Range("O245").Value = "Provincia de Cautín"
$ws.Range("P245").Value = 583
$ws.Range("Q245").Value = 12
$ws.Range("R245").Value = "Hortaliza"
